$d = $word.ActiveDocument

# Locate the paragraph that contains the literal placeholder text
# "{{generationChart}}" (the template tag docxtemplater replaces with a
# generated chart image).
$chartParaIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text -like "*{{generationChart}}*") {
        $chartParaIndex = $i
        break
    }
}

if ($chartParaIndex -lt 0) {
    throw "Could not find the {{generationChart}} placeholder paragraph"
}

# The paragraph immediately before it is the blank spacer paragraph that
# was inserted right after the "- Electricity Generation:" line. Both of
# these paragraphs should be removed so that the "Electricity Generation:"
# paragraph is immediately followed by the paragraph that originally came
# after the placeholder.
$spacerParaIndex = $chartParaIndex - 1
$spacerPara = $d.Paragraphs.Item($spacerParaIndex)

if ($spacerPara.Range.Text.Trim() -ne "") {
    throw "Expected the paragraph before the {{generationChart}} placeholder to be blank"
}

$chartPara = $d.Paragraphs.Item($chartParaIndex)

# Delete the placeholder paragraph (including its own paragraph mark) first,
# then delete the now-re-indexed blank spacer paragraph (including its own
# paragraph mark). Deleting bottom-most first keeps ranges/indices valid.
$chartPara.Range.Delete()
$spacerPara.Range.Delete()
